$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1456.942
$ws.Range("I15").Value = 1456.942
$ws.Range("K15").Value = 4370.826
$ws.Range("M15").Value = -4201.826
$ws.Range("H17").Value = 2080
$ws.Range("J17").Value = 2135.652
$ws.Range("L17").Value = 6406.956
$ws.Range("N17").Value = -6742.956
$ws.Range("H74").Value = 3354
$ws.Range("I74").Value = 3315
$ws.Range("K74").Value = 3315
$ws.Range("M74").Value = -2379
$ws.Range("H77").Value = 3354
$ws.Range("I77").Value = 3315
$ws.Range("K77").Value = 16575
$ws.Range("M77").Value = -11895
$ws.Range("H100").Value = 1988.75
$ws.Range("I100").Value = 1701.4286
$ws.Range("K100").Value = 1701.4286
$ws.Range("M100").Value = -1160.4286
$ws.Range("H113").Value = 4512.615
$ws.Range("I113").Value = 3878.3333
$ws.Range("J113").Value = 4848.4116
$ws.Range("K113").Value = 3878.3333
$ws.Range("L113").Value = 4848.4116
$ws.Range("M113").Value = -624.3332999999998
$ws.Range("N113").Value = -11356.4116
$ws.Range("H127").Value = 32259306
$ws.Range("I127").Value = 142857500
$ws.Range("J127").Value = 1497.5
$ws.Range("K127").Value = 428572500
$ws.Range("L127").Value = 4492.5
$ws.Range("M127").Value = -428567540
$ws.Range("N127").Value = -14412.5
$ws.Range("H129").Value = 635.6316
$ws.Range("I129").Value = 380.86667
$ws.Range("K129").Value = 1142.60001
$ws.Range("M129").Value = 3857.39999
$ws.Range("H132").Value = 6639.0435
$ws.Range("I132").Value = 2335.6
$ws.Range("J132").Value = 35328.668
$ws.Range("K132").Value = 7006.799999999999
$ws.Range("L132").Value = 105986.004
$ws.Range("M132").Value = -4476.799999999999
$ws.Range("N132").Value = -111046.004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5686.443
$ws.Range("I32").Value = 3937.9194
$ws.Range("J32").Value = 19237.5
$ws.Range("K32").Value = 3937.9194
$ws.Range("L32").Value = 19237.5
$ws.Range("M32").Value = -3650.9194
$ws.Range("N32").Value = -19811.5
$ws.Range("H48").Value = 100000
$ws.Range("J48").Value = 100000
$ws.Range("L48").Value = 100000
$ws.Range("N48").Value = -100768
$ws.Range("H61").Value = 1951.5555
$ws.Range("I61").Value = 1893.8823
$ws.Range("J61").Value = 2129.818
$ws.Range("K61").Value = 1893.8823
$ws.Range("L61").Value = 2129.818
$ws.Range("M61").Value = -1681.8823
$ws.Range("N61").Value = -2553.818
$ws.Range("H133").Value = 33199.8
$ws.Range("J133").Value = 33199.8
$ws.Range("L133").Value = 33199.8
$ws.Range("N133").Value = -38259.8
$ws.Range("H136").Value = 1951.5555
$ws.Range("I136").Value = 1893.8823
$ws.Range("J136").Value = 2129.818
$ws.Range("K136").Value = 5681.6469
$ws.Range("L136").Value = 6389.454000000001
$ws.Range("M136").Value = -3131.6469
$ws.Range("N136").Value = -11489.454
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 102842
$ws.Range("J43").Value = 102842
$ws.Range("L43").Value = 102842
$ws.Range("N43").Value = -103204
$ws.Range("H99").Value = 2171.4285
$ws.Range("I99").Value = 1741.6666
$ws.Range("J99").Value = 2493.75
$ws.Range("K99").Value = 1741.6666
$ws.Range("L99").Value = 2493.75
$ws.Range("M99").Value = -243.6666
$ws.Range("N99").Value = -5489.75
$ws.Range("H129").Value = 38571.285
$ws.Range("J129").Value = 38571.285
$ws.Range("L129").Value = 38571.285
$ws.Range("N129").Value = -48571.285
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 211
$ws.Range("I2").Value = 211
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 211
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -98
$ws.Range("N2").ClearContents()
$ws.Range("H88").Value = 27085.75
$ws.Range("J88").Value = 27085.75
$ws.Range("L88").Value = 27085.75
$ws.Range("N88").Value = -27897.75
$ws.Range("H91").Value = 27085.75
$ws.Range("J91").Value = 27085.75
$ws.Range("L91").Value = 27085.75
$ws.Range("N91").Value = -29893.75
$ws.Range("H141").Value = 67828.2
$ws.Range("J141").Value = 67828.2
$ws.Range("L141").Value = 67828.2
$ws.Range("N141").Value = -78188.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53207.105
$ws.Range("I2").Value = 22.714285
$ws.Range("J2").Value = 84231.336
$ws.Range("K2").Value = 136.28571
$ws.Range("L2").Value = 505388.0159999999
$ws.Range("M2").Value = -23.28570999999999
$ws.Range("N2").Value = -505614.0159999999
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H113").Value = 2020737.1
$ws.Range("I113").Value = 3367503.8
$ws.Range("J113").Value = 587.1667
$ws.Range("K113").Value = 10102511.4
$ws.Range("L113").Value = 1761.5001
$ws.Range("M113").Value = -10100341.4
$ws.Range("N113").Value = -6101.5001
$ws.Range("H131").Value = 968.7263
$ws.Range("J131").Value = 983.95605
$ws.Range("L131").Value = 2951.86815
$ws.Range("N131").Value = -13031.86815
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3027.3635
$ws.Range("I82").Value = 3171.5715
$ws.Range("J82").Value = 2775
$ws.Range("K82").Value = 3171.5715
$ws.Range("L82").Value = 2775
$ws.Range("M82").Value = -2810.5715
$ws.Range("N82").Value = -3497
$ws.Range("H85").Value = 3027.3635
$ws.Range("I85").Value = 3171.5715
$ws.Range("J85").Value = 2775
$ws.Range("K85").Value = 3171.5715
$ws.Range("L85").Value = 2775
$ws.Range("M85").Value = -1923.5715
$ws.Range("N85").Value = -5271
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 11998
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 11998
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 11998
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -12450
$ws.Range("H132").Value = 1281254.8
$ws.Range("I132").Value = 1784846.1
$ws.Range("J132").Value = 2907.3076
$ws.Range("K132").Value = 5354538.300000001
$ws.Range("L132").Value = 8721.9228
$ws.Range("M132").Value = -5352008.300000001
$ws.Range("N132").Value = -13781.9228
